$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 2019 (P) and 2020 (Q) columns with revised figures, then drop
# the 2021/2022 (R/S) columns entirely - the source table now only reports
# through 2020.
$ws.Range("P5").Value = 27
$ws.Range("Q5").Value = 25.3

$ws.Range("P6").Value = 19.6
$ws.Range("Q6").Value = 17.8

$ws.Range("P8").Value = 2.2
$ws.Range("Q8").Value = 2

$ws.Range("P9").Value = 5.2
$ws.Range("Q9").Value = 5.5

# Remove columns R:S (years 2021 and 2022) for the whole table, shifting
# nothing else - these columns only had data in rows 4-10.
$ws.Range("R4:S10").Delete()

# Match the author's final selection state.
$ws.Range("N13").Select()
